# Set percent of fuel demand changes that affect exports to 1
# Replaces the diagonal (From type == To type) formulas on sheet "PoFDCtAE"
# with a flat value of 1. Cells such as R10 (=1-J10) etc. keep their existing
# formulas and will simply recalculate to 0 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoFDCtAE")

$ws.Range("C3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("N14").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("T20").Value = 1

$excel.Calculate()
